$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @{r=2; c=4; v="330.33"},
  @{r=2; c=5; v="6.19%"},
  @{r=2; c=7; v="5"},
  @{r=3; c=4; v="40.06"},
  @{r=3; c=5; v="6.26%"},
  @{r=3; c=7; v="5"},
  @{r=4; c=4; v="5.279"},
  @{r=4; c=5; v="2.21%"},
  @{r=4; c=7; v="5"},
  @{r=5; c=4; v="0.08109"},
  @{r=5; c=5; v="2.23%"},
  @{r=5; c=7; v="5"},
  @{r=6; c=2; v="KuCoinToken"},
  @{r=6; c=3; v="https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"},
  @{r=6; c=4; v="8.637"},
  @{r=6; c=5; v="4.24%"},
  @{r=6; c=7; v="5"},
  @{r=7; c=2; v="FTXToken"},
  @{r=7; c=3; v="https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"},
  @{r=7; c=4; v="1.915"},
  @{r=7; c=5; v="-0.32%"},
  @{r=7; c=7; v="5"},
  @{r=8; c=2; v="BTSEToken"},
  @{r=8; c=3; v="https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"},
  @{r=8; c=4; v="2.959"},
  @{r=8; c=5; v="-1.00%"},
  @{r=8; c=7; v="5"},
  @{r=9; c=2; v="MXToken"},
  @{r=9; c=3; v="https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"},
  @{r=9; c=4; v="0.9365"},
  @{r=9; c=5; v="0.57%"},
  @{r=9; c=7; v="5"},
  @{r=10; c=2; v="LiechtensteinCryptoassetsExchange"},
  @{r=10; c=3; v="https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"},
  @{r=10; c=4; v="0.1377"},
  @{r=10; c=5; v="26.46%"},
  @{r=10; c=7; v="5"},
  @{r=11; c=2; v="WazirX"},
  @{r=11; c=3; v="https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"},
  @{r=11; c=4; v="0.1963"},
  @{r=11; c=5; v="1.91%"},
  @{r=11; c=7; v="5"},
  @{r=12; c=2; v="MandalaExchangeToken"},
  @{r=12; c=3; v="https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"},
  @{r=12; c=4; v="0.09197"},
  @{r=12; c=5; v="0.43%"},
  @{r=12; c=7; v="5"},
  @{r=13; c=2; v="BitrueCoin"},
  @{r=13; c=3; v="https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"},
  @{r=13; c=4; v="0.03441"},
  @{r=13; c=5; v="4.34%"},
  @{r=13; c=7; v="5"},
  @{r=14; c=2; v="BitMartToken"},
  @{r=14; c=3; v="https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"},
  @{r=14; c=4; v="0.09560"},
  @{r=14; c=5; v="-0.41%"},
  @{r=14; c=7; v="5"},
  @{r=15; c=2; v="BitForexToken"},
  @{r=15; c=3; v="https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"},
  @{r=15; c=4; v="0.001395"},
  @{r=15; c=5; v="1.26%"},
  @{r=15; c=7; v="5"},
  @{r=16; c=2; v="TigerCash"},
  @{r=16; c=3; v="https://coinranking.com/coin/6hIn06L2+tigercash-tch"},
  @{r=16; c=4; v="0.005898"},
  @{r=16; c=5; v="1.42%"},
  @{r=16; c=7; v="5"},
  @{r=17; c=2; v="LEO"},
  @{r=17; c=3; v="https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"},
  @{r=17; c=4; v="3.361"},
  @{r=17; c=5; v="-6.52%"},
  @{r=17; c=7; v="5"},
  @{r=18; c=2; v="GateToken"},
  @{r=18; c=3; v="https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"},
  @{r=18; c=4; v="4.530"},
  @{r=18; c=5; v="2.05%"},
  @{r=18; c=7; v="5"},
  @{r=19; c=4; v="0.3524"},
  @{r=19; c=5; v="3.36%"},
  @{r=19; c=7; v="5"},
  @{r=20; c=4; v="7.422"},
  @{r=20; c=5; v="15.62%"},
  @{r=20; c=7; v="5"},
  @{r=21; c=5; v="1.97%"},
  @{r=21; c=7; v="5"},
  @{r=22; c=4; v="0.2311"},
  @{r=22; c=5; v="-10.81%"},
  @{r=22; c=7; v="5"},
  @{r=23; c=4; v="0.04443"},
  @{r=23; c=5; v="0.76%"},
  @{r=23; c=7; v="5"},
  @{r=24; c=4; v="0.001223"},
  @{r=24; c=5; v="-0.86%"},
  @{r=24; c=7; v="5"},
  @{r=25; c=4; v="0.004357"},
  @{r=25; c=5; v="-5.84%"},
  @{r=25; c=7; v="5"},
  @{r=26; c=5; v="-5.18%"},
  @{r=26; c=7; v="5"},
  @{r=27; c=5; v="-0.07%"},
  @{r=27; c=7; v="5"},
  @{r=28; c=7; v="5"},
  @{r=29; c=7; v="5"},
  @{r=30; c=7; v="5"},
  @{r=31; c=7; v="5"},
  @{r=32; c=7; v="5"},
  @{r=33; c=7; v="5"},
  @{r=34; c=7; v="5"},
  @{r=35; c=7; v="5"},
  @{r=36; c=7; v="5"},
  @{r=37; c=7; v="5"},
  @{r=38; c=7; v="5"},
  @{r=39; c=4; v="0.02511"},
  @{r=39; c=5; v="11.86%"},
  @{r=39; c=7; v="5"},
  @{r=40; c=4; v="0.05233"},
  @{r=40; c=5; v="2.82%"},
  @{r=40; c=7; v="5"},
  @{r=41; c=4; v="0.007688"},
  @{r=41; c=5; v="2.97%"},
  @{r=41; c=7; v="5"},
  @{r=42; c=4; v="0.1430"},
  @{r=42; c=5; v="5.45%"},
  @{r=42; c=7; v="5"},
  @{r=43; c=4; v="0.009020"},
  @{r=43; c=5; v="0.89%"},
  @{r=43; c=7; v="5"},
  @{r=44; c=4; v="0.002170"},
  @{r=44; c=5; v="1.85%"},
  @{r=44; c=7; v="5"},
  @{r=45; c=4; v="0.008133"},
  @{r=45; c=7; v="5"},
  @{r=46; c=4; v="0.00006660"},
  @{r=46; c=5; v="0.55%"},
  @{r=46; c=7; v="5"},
  @{r=47; c=5; v="-0.08%"},
  @{r=47; c=7; v="5"},
  @{r=48; c=2; v="BOLO"},
  @{r=48; c=3; v="https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"},
  @{r=48; c=4; v="0.003337"},
  @{r=48; c=5; v="16.55%"},
  @{r=48; c=7; v="5"},
  @{r=49; c=2; v="CoinbaseStockToken"},
  @{r=49; c=3; v="https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"},
  @{r=49; c=4; v="0.002483"},
  @{r=49; c=5; v="148.20%"},
  @{r=49; c=7; v="5"},
  @{r=50; c=5; v="-0.08%"},
  @{r=50; c=7; v="5"},
  @{r=51; c=5; v="-0.08%"},
  @{r=51; c=7; v="5"}
)

foreach ($item in $data) {
    $cell = $ws.Cells.Item($item.r, $item.c)
    $cell.NumberFormat = "@"
    $cell.Value = $item.v
}
